$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.590.72'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '2.417.15'
$ws.Range('E3').Value = '  +2.07%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.68'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.43'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.508'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.01'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0796'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('E12').Value = '  +2.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.50'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.89'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').Value = '2.788.58'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').Value = '2.424.42'
$ws.Range('E16').Value = '  +2.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.827'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.43%  '
$ws.Range('D18').Value = '43.641.42'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.41'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.08'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.24%  '
$ws.Range('D21').Value = '0.0₃0899'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.15'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.31'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.25'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.94'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.21'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.49%  '
$ws.Range('E29').Value = '  +2.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.33'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.120'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +18.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.13'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.38'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.36%  '
$ws.Range('E34').Value = '  +0.15%  '
$ws.Range('E35').Value = '  +2.81%  '
$ws.Range('E36').Value = '  +2.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '131.13'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +27.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.91'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +4.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.38'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.108'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.17'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.78%  '
$ws.Range('D43').Value = '1.943.39'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0283'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.17%  '
$ws.Range('E45').Value = '  +2.29%  '
$ws.Range('E46').Value = '  +2.89%  '
$ws.Range('E47').Value = '  -1.83%  '
$ws.Range('D48').Value = '2.647.36'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('E49').Value = '  +3.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.69'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.27'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.11%  '
